$d = $word.ActiveDocument

# --- Change 1: "API: Monads / Transforms" -> "Functional API: Monads / Transforms" ---
$rng1 = $d.Content
$rng1.Find.Execute("API: Monads / Transforms", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng1.Collapse(1) | Out-Null
$rng1.InsertBefore("Functional ")

# --- Change 2: add two new paragraphs right after
#     "Mapping: Selector Monad. Matching Resource / Role set?" ---
$rng2 = $d.Content
$rng2.Find.Execute("Mapping: Selector Monad. Matching Resource / Role set?", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorIndex = $rng2.Paragraphs(1).Index

# new blank paragraph right after the anchor
$d.Paragraphs($anchorIndex).Range.InsertParagraphAfter() | Out-Null

# new paragraph with the added text, right after the blank one
$d.Paragraphs($anchorIndex + 1).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs($anchorIndex + 2).Range.InsertBefore("Context / Occurrence Monads Instances (Layers Hierarchies Monads): Aligned / Matching Entities resolution (Augmentations Agggregations / Activation / Alignments matchings). Versioned graph: stateless / functional. Mappings assertions matching.")
